$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 2; $i -le 480; $i++) {
    $ws.Cells.Item($i, 3).Value = 45177
}

$ws.Rows.Item(480).RowHeight = 15

$ws.Cells.Item(481, 1).Value = "A 41729-2023"
$ws.Cells.Item(481, 2).Value = 45176
$ws.Cells.Item(481, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(481, 3).Value = 45177
$ws.Cells.Item(481, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(481, 4).Value = "GÄVLEBORGS LÄN"
$ws.Cells.Item(481, 5).Value = "NORDANSTIG"
$ws.Cells.Item(481, 6).Value = "Holmen skog AB"
$ws.Cells.Item(481, 7).Value = 1
$ws.Cells.Item(481, 8).Value = 0
$ws.Cells.Item(481, 9).Value = 0
$ws.Cells.Item(481, 10).Value = 0
$ws.Cells.Item(481, 11).Value = 0
$ws.Cells.Item(481, 12).Value = 0
$ws.Cells.Item(481, 13).Value = 0
$ws.Cells.Item(481, 14).Value = 0
$ws.Cells.Item(481, 15).Value = 0
$ws.Cells.Item(481, 16).Value = 0
$ws.Cells.Item(481, 17).Value = 0
$ws.Cells.Item(481, 18).WrapText = $true
$ws.Rows.Item(481).RowHeight = 15

$ws.Cells.Item(482, 1).Value = "A 41721-2023"
$ws.Cells.Item(482, 2).Value = 45176
$ws.Cells.Item(482, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(482, 3).Value = 45177
$ws.Cells.Item(482, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(482, 4).Value = "GÄVLEBORGS LÄN"
$ws.Cells.Item(482, 5).Value = "NORDANSTIG"
$ws.Cells.Item(482, 7).Value = 8.1
$ws.Cells.Item(482, 8).Value = 0
$ws.Cells.Item(482, 9).Value = 0
$ws.Cells.Item(482, 10).Value = 0
$ws.Cells.Item(482, 11).Value = 0
$ws.Cells.Item(482, 12).Value = 0
$ws.Cells.Item(482, 13).Value = 0
$ws.Cells.Item(482, 14).Value = 0
$ws.Cells.Item(482, 15).Value = 0
$ws.Cells.Item(482, 16).Value = 0
$ws.Cells.Item(482, 17).Value = 0
$ws.Cells.Item(482, 18).WrapText = $true
